# Daily attendance processing - 2025-12-31 06:03:15
#
# Reorders the "Recorded By" values in column G so that the literal
# entry "System" (exact case) is moved to the front of the
# comma-separated list, keeping the relative order of the remaining
# entries (e.g. "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com").
#
# Note: PowerShell's case-sensitive operators (-ceq/-cne/-cmatch) are not
# reliable in this host, so exact-case comparisons use the .NET
# [string]::Equals / .Equals() method instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($null -eq $val) { continue }

    $text = [string]$val
    if ($text -eq "") { continue }

    $parts = $text -split ", "
    if ($parts.Count -le 1) { continue }

    $hasSystem = $false
    foreach ($p in $parts) {
        if ($p.Equals("System")) { $hasSystem = $true }
    }
    if (-not $hasSystem) { continue }

    $rest = @()
    $alreadyRemovedOne = $false
    foreach ($p in $parts) {
        if ((-not $alreadyRemovedOne) -and $p.Equals("System")) {
            $alreadyRemovedOne = $true
        } else {
            $rest += $p
        }
    }

    $newParts = @("System") + $rest
    $newText = $newParts -join ", "

    if (-not $newText.Equals($text)) {
        $cell.Value2 = $newText
    }
}
